$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2:F25 block (columns B, C, D, E, F)
$dataLeft = @(
    @(1.02,1.018022773874987,1.02363140776266,1.019342340285812,1.032005233336023),
    @(1.02,1.018923337957142,1.024279362636225,1.020104939593014,1.034123244324402),
    @(1.02,1.019506446566109,1.024698666396276,1.020599111908618,1.035488244591013),
    @(1.02,1.019751676810899,1.024874949268084,1.020807033361278,1.036060805471205),
    @(1.02,1.019792857394907,1.024904548367578,1.020841954314596,1.036156866464221),
    @(1.02,1.019509722986091,1.024701021867536,1.020601889492818,1.03549590019093),
    @(1.02,1.018327044089732,1.023850379900607,1.019599914790248,1.032722181001849),
    @(1.02,1.016245966409975,1.022351717458747,1.017839855292733,1.027791086241886),
    @(1.02,1.014860586902881,1.021352819124251,1.016670259890304,1.024472551623224),
    @(1.02,1.014261181118672,1.020920338830324,1.01616471745988,1.023027770426932),
    @(1.02,1.014038606282244,1.020759704011367,1.015977072406366,1.02248989977004),
    @(1.02,1.0140863461179,1.020794160377416,1.016017316746106,1.022605330425918),
    @(1.02,1.014242781560833,1.020907060549526,1.016149203884715,1.022983334868102),
    @(1.02,1.014339176033733,1.020976623087192,1.016230481933781,1.023216073797899),
    @(1.02,1.014900377462347,1.021381522474232,1.016703830099359,1.024568268744666),
    @(1.02,1.015252530949751,1.021635518515573,1.017000990368526,1.025414341200626),
    @(1.02,1.015457981720653,1.021783674851828,1.017174405728297,1.025907087676084),
    @(1.02,1.015528042853551,1.021834193130854,1.017233550586223,1.026074974852069),
    @(1.02,1.015214743495838,1.021608266652104,1.016969098900741,1.025323643818534),
    @(1.02,1.014196713256483,1.020873814050909,1.016110362652282,1.022872055815644),
    @(1.02,1.013557049010641,1.020412078713231,1.015571227773718,1.021323598758105),
    @(1.02,1.013896108054482,1.020656849101142,1.01585695858086,1.022145145797456),
    @(1.02,1.015231817879331,1.021620580581896,1.016983509014526,1.02536462838764),
    @(1.02,1.016783622762638,1.022739122239997,1.018294210448081,1.029071237582558)
)

# I2:N25 block (columns I, J, K, L, M, N)
$dataRight = @(
    @(1.028446458462967,1.023233566503164,1.026462277794471,1.022185871860443,1.034811713168608,1.02468667553413),
    @(1.028627774532127,1.023770528337955,1.026917585974439,1.022754594950167,1.03673490464504,1.025224399916316),
    @(1.028743417867088,1.024117760876022,1.027211551032819,1.023122658112442,1.037973504155172,1.025572125564415),
    @(1.028791633196816,1.024263684993294,1.02733497832473,1.023277406033094,1.038492842756571,1.025718256910672),
    @(1.028799705273007,1.024288183209638,1.027355693182514,1.023303389705714,1.038579962542896,1.025742789917292),
    @(1.028744063697477,1.024119710927721,1.027213200885925,1.023124725806726,1.037980448926717,1.025574078385411),
    @(1.028508083682836,1.023415080146007,1.026616285836355,1.022378061499326,1.03546288968036,1.024868446947165),
    @(1.028079335605302,1.022171777163913,1.025559470077662,1.021062834544836,1.030980677161311,1.023623378332228),
    @(1.027784740812327,1.021341812561962,1.024851579807311,1.020186373416403,1.027959848401208,1.022792235085344),
    @(1.027655082125253,1.020982171290274,1.024544260443346,1.019806945633585,1.026643639923199,1.022432083081798),
    @(1.027606604450699,1.020848545394667,1.024429988209711,1.019666022456809,1.026153479038948,1.022298267422092),
    @(1.027617017432054,1.020877210371208,1.024454505432,1.019696250354429,1.02625867791015,1.022326973106188),
    @(1.027651081414727,1.020971126529373,1.024534817117003,1.019795296608826,1.026603149030171,1.02242102263607),
    @(1.027672027349448,1.021028986194181,1.024584283870308,1.019856324033088,1.026815220710173,1.022478964468239),
    @(1.027793301512391,1.021365675277528,1.024871958746135,1.020211556621819,1.028047025719338,1.022816131688703),
    @(1.02786881100408,1.021576801622126,1.025052195747516,1.020434407810202,1.028817492647139,1.023027557856926),
    @(1.027912652173242,1.02169992284653,1.025157247942644,1.020564401383881,1.029266107111518,1.023150853927591),
    @(1.027927566622843,1.021741899723731,1.025193054971995,1.020608727221348,1.029418940806923,1.023192890416773),
    @(1.027860730480272,1.021554152369273,1.025032865988944,1.020410497131569,1.028734910374467,1.023004876439537),
    @(1.027641059175272,1.020943471640493,1.024511170645762,1.019766129588775,1.026501746005341,1.022393328474076),
    @(1.027501110167974,1.020559286169079,1.024182464867103,1.019361066634459,1.025090344974615,1.022008597415226),
    @(1.027575474064859,1.020762971568942,1.024356783996384,1.019575790815462,1.025839261262752,1.022212572071719),
    @(1.027864382346124,1.021564386669906,1.025041600512997,1.02042130131289,1.028772228175726,1.023015125274051),
    @(1.028191716118515,1.0224933950495,1.025833271928691,1.021402791244801,1.03214506469944,1.023945452952102)
)

for ($i = 0; $i -lt $dataLeft.Count; $i++) {
    $row = $dataLeft[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $excelCol = $j + 2   # column B = 2
        $ws.Cells.Item($excelRow, $excelCol).Value = $row[$j]
    }
}

for ($i = 0; $i -lt $dataRight.Count; $i++) {
    $row = $dataRight[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $excelCol = $j + 9   # column I = 9
        $ws.Cells.Item($excelRow, $excelCol).Value = $row[$j]
    }
}
